$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Student data (rows 3-9) ------------------------------------------------
# Column A holds numeric-looking student IDs that must stay TEXT, not numbers.
# A leading apostrophe forces text entry (quotePrefix); resetting the style
# back to "Normal" afterwards drops the cell back onto the sheet's default
# style, matching a cell that was always typed/imported as plain text.
#
# Column A is filled top-to-bottom first, then columns B/C are filled in per
# row - this mirrors how the workbook's shared-strings table ends up ordered
# (all IDs first, then each row's first/last name pair) after the edit.

$ids = "6110620013", "6110620033", "6110620034", "6110620035", "6110620038", "6110620045", "6110620048"
for ($i = 0; $i -lt $ids.Length; $i++) {
    $cell = $ws.Cells.Item(3 + $i, 1)
    $cell.Value = "'" + $ids[$i]
    $cell.Style = "Normal"
}

$names = @(
    @("น.ส. ฐิตาภรณ์ ", "อนุสาร"),
    @("นาย กิตตินันท์ ", "ขวัญซ้าย"),
    @("น.ส. เขมินี ", "ทองมา"),
    @("น.ส. จอมทอง ", "ชัยภักดี"),
    @("น.ส. ทิพยเนตร ", "จงจิตร"),
    @("นาย พีระพัฒน์ ", "เพ่งพิศ"),
    @("น.ส. สุชานาถ ", "กองวารี")
)
for ($i = 0; $i -lt $names.Length; $i++) {
    $r = 3 + $i
    $ws.Cells.Item($r, 2).Value = $names[$i][0]
    $ws.Cells.Item($r, 3).Value = $names[$i][1]
}

# Row 6 used to carry a note in column D ("โทร.0899454565464") - remove it.
$ws.Cells.Item(6, 4).ClearContents()

# --- Column widths (closest reachable match to the new layout; the host
#     quantizes ColumnWidth to 1/7-character steps, so these are the inputs
#     that land nearest the target stored widths of 11.875/21.75/22.25/17.625)
$ws.Columns.Item(1).ColumnWidth = 11.142857142857142
$ws.Columns.Item(2).ColumnWidth = 21.0
$ws.Columns.Item(3).ColumnWidth = 21.571428571428573
$ws.Columns.Item(4).ColumnWidth = 16.857142857142858

# --- Selection ---------------------------------------------------------
$ws.Range("C9").Select() | Out-Null
